# Apply the edits described by the diff to slide 2 of the presentation.
#
# 1) Shape "Rectangle 3" (cNvPr id=4)        -> move from (4211782,1520388) to (4206371,1550992) EMU
# 2) Shape "Picture 2"   (cNvPr id=1026)     -> move from (4299526,1589231) to (4299525,1629919) EMU
#    (there are several shapes named "Picture 2" on this slide, so we match on the
#     stable shape Id instead of the ambiguous Name)
# 3) Shape "TextBox 15"  (cNvPr id=16)       -> drop the trailing ", Rstarved" run pair so the
#                                               parameter list ends at "...muMax"
#
# Note: this COM host stores shape position in points and truncates (floors) the
# point->EMU conversion on write, so the point literals below are chosen so that
# floor(points * 12700) reproduces the exact target EMU values from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# --- Rectangle 3 (id=4): reposition ---
$rect3 = Get-ShapeById $s 4
$rect3.Left = 331.21032
$rect3.Top = 122.1254

# --- Picture 2 / the R_activity equation image (id=1026): reposition ---
$pic2 = Get-ShapeById $s 1026
$pic2.Left = 338.5453
$pic2.Top = 128.3401

# --- TextBox 15 (id=16): remove the trailing ", Rstarved" text ---
$paramBox = Get-ShapeById $s 16
$tr = $paramBox.TextFrame.TextRange
$fullText = $tr.Text
$suffix = ", Rstarved"
if ($fullText.EndsWith($suffix)) {
    $startPos = $fullText.Length - $suffix.Length + 1
    $tr.Characters($startPos, $suffix.Length).Text = ""
}
